# The deck's "Insert Date and Time" placeholder (a fixed, non auto-updating
# "datetime1" field) was refreshed from 3/30/2020 to 4/3/2020. That
# placeholder lives on the slide master and on every slide layout (no
# individual slide overrides it), so walk Master + all CustomLayouts and
# update the one shape whose PlaceholderFormat.Type is the date placeholder
# (ppPlaceholderDate = 16).

$p = $ppt.ActivePresentation
$newDate = "4/3/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shp.HasTextFrame) {
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout has its own (inherited-looking, but separately stored)
# date placeholder shape as well.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
